$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptos list values (price + volume%) as of the latest scrape.
# Leading apostrophe forces Excel to store the value as literal text,
# matching the original inline-string cell type instead of auto-converting
# numeric-looking strings (e.g. '568.26') into real numbers.

$ws.Range("D2").Value = "'61.212.20"
$ws.Range("E2").Value = "'  -1.84%  "
$ws.Range("D3").Value = "'3.013.29"
$ws.Range("E3").Value = "'  -3.87%  "
$ws.Range("E4").Value = "'  -0.14%  "
$ws.Range("D5").Value = "'568.26"
$ws.Range("E5").Value = "'  -3.03%  "
$ws.Range("D6").Value = "'128.78"
$ws.Range("E6").Value = "'  -4.58%  "
$ws.Range("E7").Value = "'  -0.16%  "
$ws.Range("D8").Value = "'3.012.63"
$ws.Range("E8").Value = "'  -3.62%  "
$ws.Range("E9").Value = "'  -1.70%  "
$ws.Range("D10").Value = "'0.134"
$ws.Range("E10").Value = "'  -4.34%  "
$ws.Range("E11").Value = "'  -0.83%  "
$ws.Range("E12").Value = "'  -4.70%  "
$ws.Range("E13").Value = "'  -4.46%  "
$ws.Range("D14").Value = "'32.77"
$ws.Range("E14").Value = "'  -3.13%  "
$ws.Range("E15").Value = "'  +0.10%  "
$ws.Range("D16").Value = "'3.510.55"
$ws.Range("E16").Value = "'  -3.86%  "
$ws.Range("D17").Value = "'61.252.72"
$ws.Range("E17").Value = "'  -1.74%  "
$ws.Range("D18").Value = "'3.014.21"
$ws.Range("E18").Value = "'  -3.82%  "
$ws.Range("D19").Value = "'6.22"
$ws.Range("E19").Value = "'  -4.49%  "
$ws.Range("D20").Value = "'438.23"
$ws.Range("E20").Value = "'  -2.90%  "
$ws.Range("D21").Value = "'13.15"
$ws.Range("E21").Value = "'  -4.77%  "
$ws.Range("E22").Value = "'  -5.39%  "
$ws.Range("E23").Value = "'  -5.28%  "
$ws.Range("D24").Value = "'78.98"
$ws.Range("E24").Value = "'  -4.95%  "
$ws.Range("E25").Value = "'  -5.55%  "
$ws.Range("E26").Value = "'  -0.07%  "
$ws.Range("D27").Value = "'1.00"
$ws.Range("E27").Value = "'  -0.10%  "
$ws.Range("E28").Value = "'  -6.42%  "
$ws.Range("E29").Value = "'  -5.93%  "
$ws.Range("E30").Value = "'  -7.41%  "
$ws.Range("D31").Value = "'25.58"
$ws.Range("E31").Value = "'  -5.32%  "
$ws.Range("E32").Value = "'  -6.22%  "
$ws.Range("E33").Value = "'  -8.28%  "
$ws.Range("D34").Value = "'2.26"
$ws.Range("E34").Value = "'  -3.79%  "
$ws.Range("D35").Value = "'0.955"
$ws.Range("E35").Value = "'  -6.49%  "
$ws.Range("D36").Value = "'5.54"
$ws.Range("E36").Value = "'  -3.74%  "
$ws.Range("D37").Value = "'50.08"
$ws.Range("E37").Value = "'  -1.82%  "
$ws.Range("D38").Value = "'0.0₃0674"
$ws.Range("E38").Value = "'  -3.16%  "
$ws.Range("E39").Value = "'  -5.13%  "
$ws.Range("B40").Value = "'Kaspa"
$ws.Range("C40").Value = "'https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D40").Value = "'0.109"
$ws.Range("E40").Value = "'  -1.66%  "
$ws.Range("B41").Value = "'Cosmos"
$ws.Range("C41").Value = "'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D41").Value = "'7.72"
$ws.Range("E41").Value = "'  -3.40%  "
$ws.Range("D42").Value = "'373.22"
$ws.Range("E42").Value = "'  -4.70%  "
$ws.Range("D43").Value = "'2.648.15"
$ws.Range("E43").Value = "'  -3.42%  "
$ws.Range("E44").Value = "'  -8.85%  "
$ws.Range("E46").Value = "'  -4.96%  "
$ws.Range("D47").Value = "'119.68"
$ws.Range("E47").Value = "'  -3.97%  "
$ws.Range("E48").Value = "'  -6.41%  "
$ws.Range("D49").Value = "'32.91"
$ws.Range("E49").Value = "'  -3.57%  "
$ws.Range("D50").Value = "'0.107"
$ws.Range("E50").Value = "'  -3.15%  "
$ws.Range("D51").Value = "'23.56"
$ws.Range("E51").Value = "'  -5.83%  "
